$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Play Divine Showdown Free - Exciting Features & High Jackpot", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Divine Showdown for Free - Exciting Special Features!", 2)

$d.Content.Find.Execute(
    "Exciting special features that boost winnings", $true, $false, $false, $false, $false,
    $true, 1, $false, "Exciting special features", 2)

$d.Content.Find.Execute(
    "Maximum prize of €50,000", $true, $false, $false, $false, $false,
    $true, 1, $false, "High maximum prize", 2)

$d.Content.Find.Execute(
    "High-quality graphics and animations", $true, $false, $false, $false, $false,
    $true, 1, $false, "Variable RTP percentage", 2)

$d.Content.Find.Execute(
    "Developed by a reputable game developer", $true, $false, $false, $false, $false,
    $true, 1, $false, "Developed by reputable company", 2)

$d.Content.Find.Execute(
    "High volatility may not appeal to some players", $true, $false, $false, $false, $false,
    $true, 1, $false, "High volatility", 2)

$d.Content.Find.Execute(
    "Variable RTP may not be ideal for all players", $true, $false, $false, $false, $false,
    $true, 1, $false, "Betting range may not suit all players", 2)

$d.Content.Find.Execute(
    "Looking to play Divine Showdown? Read our review of this exciting game featuring special features, a high jackpot, and developed by Play 'N Go. Play for free.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Read our review of Divine Showdown, a slot game with exciting special features. Play for free!", 2)
